$d = $word.ActiveDocument

# --- Edit 1 ---------------------------------------------------------------
# Insert a new paragraph "First_time = True" right after the paragraph that
# reads "Extrakcia<br>for_prediction = False" (and right before the empty
# paragraph that follows it).
$count = $d.Paragraphs.Count
$targetIndex = 0
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Extrakcia`vfor_prediction = False`r") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    # The paragraph right after the target is the existing empty paragraph;
    # insert the new text immediately before it so the empty paragraph is
    # preserved untouched.
    $nextPara = $d.Paragraphs.Item($targetIndex + 1)
    $insertRng = $nextPara.Range
    $insertRng.Collapse(1)
    $insertRng.InsertBefore("First_time = True`r")
}

# --- Edit 2 ---------------------------------------------------------------
# Insert a new paragraph "first_time = True" right after the very last
# paragraph of the body ("for_prediction = True"), before the sectPr.
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$lastRng = $lastPara.Range
$lastRng.InsertParagraphAfter()
$newLastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newLastPara.Range.Text = "first_time = True"
